$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45207 -> 45208, i.e. 2023-10-08 -> 2023-10-09) for every data row (2-15).
foreach ($r in 2..15) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}
